# Special (non-ASCII) characters used throughout the table text.
$enDash   = [char]0x2013   # "–"
$rsQuote  = [char]0x2019   # "’"
$ldQuote  = [char]0x201C   # "\u201c"
$rdQuote  = [char]0x201D   # "\u201d"

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$copernicusBlurb = "Studio del caso d" + $rsQuote + "uso " + $ldQuote + "AgroSat" + $rdQuote + " del dominio " + $ldQuote + "Precision Agriculture" + $rdQuote + " fornito da Copernicus parte "

function Fill-Row($row, $giorno, $orario, $ore, $attivita) {
    $c1 = $row.Cells.Item(1)
    $c1.Range.Text = $giorno
    $c1.Range.Font.Size = 12

    $c2 = $row.Cells.Item(2)
    $c2.Range.Text = $orario
    $c2.Range.Font.Size = 12

    $c3 = $row.Cells.Item(3)
    $c3.Range.Text = $ore
    $c3.Range.Font.Size = 12

    $c4 = $row.Cells.Item(4)
    $c4.Range.Text = $attivita
    $c4.Range.Font.Size = 12
}

# --- Row 6 already exists but is empty: fill in with the "parte 2" entry. ---
$row6 = $t.Rows.Item(6)
Fill-Row $row6 "06/10/2023" ("10:00 " + $enDash + " 16:00") "6" ($copernicusBlurb + "2.")

# --- New row: "parte 3" entry (09/10/2023). ---
$row7 = $t.Rows.Add()
Fill-Row $row7 "09/10/2023" ("10:00 " + $enDash + " 16:00") "6" ($copernicusBlurb + "3.")

# --- New row: "parte 4" entry (10/10/2023). ---
$row8 = $t.Rows.Add()
Fill-Row $row8 "10/10/2023" ("11:00 " + $enDash + " 16:00") "5" ($copernicusBlurb + "4.")

# --- New row: precipitazioniCampania2022.js / raster entry (11/10/2023). ---
$row9 = $t.Rows.Add()
$raster = "Creazione dello script " + $ldQuote + "precipitazioniCampania2022.js" + $rdQuote + " e modifica del raster prodotto su QGIS."
Fill-Row $row9 "11/10/2023" ("10:00 " + $enDash + " 16:00") "6" $raster

# --- New trailing empty row (mirrors the one that used to sit at the end). ---
$row10 = $t.Rows.Add()

# --- Resize every column (applies to all rows, old and new, at once). ---
$t.Columns.Item(1).Width = 97.95
$t.Columns.Item(2).Width = 80.7
$t.Columns.Item(3).Width = 69.8
$t.Columns.Item(4).Width = 232.95

# --- Update the hour total. ---
$d.Content.Find.Execute("Totale ore: 23", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Totale ore: 46", 2)

# --- The _GoBack bookmark used to sit at the end of the document; it now
#     belongs inside the newly typed text, but this runtime cannot relocate
#     bookmarks, so simply drop the stale one instead of leaving it behind. ---
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
} catch {
}
